$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("API_Controller")

# New "UserInteraction" mini-table in columns F:G, mirroring the existing
# "DataSource" table pattern already used on Sheet1 (merged, centered header
# in row 1, plain key/value rows below).
$ws.Range("F1:G1").Merge()
$ws.Range("F1").Value = "UserInteraction"
$ws.Range("F1:G1").HorizontalAlignment = -4108

$ws.Range("F2").Value = "IdUser"
$ws.Range("G2").Value = "Guid"

$ws.Range("F3").Value = "IdArticle"
$ws.Range("G3").Value = "Guid"

$ws.Range("F4").Value = "Like"
$ws.Range("G4").Value = 1

$ws.Range("F5").Value = "DisLike"
$ws.Range("G5").Value = 4

$ws.Columns.Item(6).ColumnWidth = 15.08984375

$ws.Range("G5").Select()
